$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

$cells = @("F3", "F4", "F5", "F8", "F9", "F10", "F13", "F14")
foreach ($cell in $cells) {
    $ws.Range($cell).ClearContents()
}
